$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, border, centered) from G1 to H1, then set header text
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Resumen"

$ws.Range("H2").Value = "['Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio']"
$ws.Range("H3").Value = "['Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio']"
$ws.Range("H4").Value = "['Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio']"
$ws.Range("H5").Value = "['Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio']"
$ws.Range("H6").Value = "['Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio']"
$ws.Range("H7").Value = "['Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio']"
$ws.Range("H8").Value = "['Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio', 'Medio Ambiente, Medio Ambiente Ambiente Medio Medio Medio Ambiente Medio Ambiente ambiente Medio Medio ambiente Medio Ambiente Natural Medio Medio Natural Medio Ambiente medio Ambiente Medio ambiente ambiente Medio ambiente Ambiente Medio espacio Medio Ambiente. Medio Ambiente Eco Medio Ambiente espacio Medio Medio medio Medio Ambiente que Medio Ambiente Marina Medio Medio Eco Medio Medio espacio Ambiente Medio medio ambiente Medio']"
